$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text before writing so numeric-looking
# strings (e.g. "1.00", "0.999") are kept as literal text instead of
# being auto-coerced to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "57.601.36"
$ws.Cells.Item(2, 5).Value = "  -2.36%  "

$ws.Cells.Item(3, 4).Value = "2.426.47"
$ws.Cells.Item(3, 5).Value = "  -3.38%  "

$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.18%  "

$ws.Cells.Item(5, 4).Value = "514.29"
$ws.Cells.Item(5, 5).Value = "  -3.68%  "

$ws.Cells.Item(6, 4).Value = "131.84"
$ws.Cells.Item(6, 5).Value = "  -2.84%  "

$ws.Cells.Item(7, 4).Value = "0.998"
$ws.Cells.Item(7, 5).Value = "  -0.15%  "

$ws.Cells.Item(8, 4).Value = "0.551"
$ws.Cells.Item(8, 5).Value = "  -2.95%  "

$ws.Cells.Item(9, 4).Value = "2.427.42"
$ws.Cells.Item(9, 5).Value = "  -3.55%  "

$ws.Cells.Item(10, 4).Value = "0.0956"
$ws.Cells.Item(10, 5).Value = "  -5.73%  "

$ws.Cells.Item(11, 5).Value = "  -1.55%  "

$ws.Cells.Item(12, 4).Value = "5.21"
$ws.Cells.Item(12, 5).Value = "  -4.05%  "

$ws.Cells.Item(13, 4).Value = "0.332"
$ws.Cells.Item(13, 5).Value = "  -4.00%  "

$ws.Cells.Item(14, 4).Value = "2.856.95"
$ws.Cells.Item(14, 5).Value = "  -3.36%  "

$ws.Cells.Item(15, 4).Value = "57.513.96"
$ws.Cells.Item(15, 5).Value = "  -2.36%  "

$ws.Cells.Item(16, 5).Value = "  -5.61%  "

$ws.Cells.Item(17, 4).Value = "0.0000133"
$ws.Cells.Item(17, 5).Value = "  -3.96%  "

$ws.Cells.Item(18, 4).Value = "2.425.02"
$ws.Cells.Item(18, 5).Value = "  -3.23%  "

$ws.Cells.Item(19, 4).Value = "10.41"
$ws.Cells.Item(19, 5).Value = "  -5.69%  "

$ws.Cells.Item(20, 4).Value = "315.38"
$ws.Cells.Item(20, 5).Value = "  -2.33%  "

$ws.Cells.Item(21, 4).Value = "4.10"
$ws.Cells.Item(21, 5).Value = "  -3.62%  "

$ws.Cells.Item(22, 4).Value = "1.00"
$ws.Cells.Item(22, 5).Value = "  -0.06%  "

$ws.Cells.Item(23, 4).Value = "5.63"
$ws.Cells.Item(23, 5).Value = "  -4.86%  "

$ws.Cells.Item(24, 4).Value = "63.66"
$ws.Cells.Item(24, 5).Value = "  -2.22%  "

$ws.Cells.Item(25, 4).Value = "0.403"
$ws.Cells.Item(25, 5).Value = "  -3.98%  "

$ws.Cells.Item(26, 4).Value = "0.998"
$ws.Cells.Item(26, 5).Value = "  -0.07%  "

$ws.Cells.Item(27, 4).Value = "0.160"
$ws.Cells.Item(27, 5).Value = "  -2.46%  "

$ws.Cells.Item(28, 4).Value = "7.22"
$ws.Cells.Item(28, 5).Value = "  -4.58%  "

$ws.Cells.Item(29, 4).Value = "169.35"
$ws.Cells.Item(29, 5).Value = "  -0.28%  "

$ws.Cells.Item(30, 4).Value = "0.0₃0723"
$ws.Cells.Item(30, 5).Value = "  -5.46%  "

$ws.Cells.Item(31, 2).Value = "Aptos"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(31, 4).Value = "6.21"
$ws.Cells.Item(31, 5).Value = "  -4.84%  "

$ws.Cells.Item(32, 2).Value = "PancakeSwap"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(32, 4).Value = "1.67"
$ws.Cells.Item(32, 5).Value = "  -4.64%  "

$ws.Cells.Item(33, 4).Value = "1.16"
$ws.Cells.Item(33, 5).Value = "  +2.44%  "

$ws.Cells.Item(34, 4).Value = "0.998"
$ws.Cells.Item(34, 5).Value = "  -0.04%  "

$ws.Cells.Item(35, 4).Value = "0.997"
$ws.Cells.Item(35, 5).Value = "  -0.15%  "

$ws.Cells.Item(36, 4).Value = "17.72"
$ws.Cells.Item(36, 5).Value = "  -3.76%  "

$ws.Cells.Item(37, 5).Value = "  -6.96%  "

$ws.Cells.Item(38, 4).Value = "3.88"
$ws.Cells.Item(38, 5).Value = "  -4.11%  "

$ws.Cells.Item(39, 4).Value = "36.30"
$ws.Cells.Item(39, 5).Value = "  -2.20%  "

$ws.Cells.Item(40, 4).Value = "1.44"
$ws.Cells.Item(40, 5).Value = "  -5.02%  "

$ws.Cells.Item(41, 4).Value = "0.776"
$ws.Cells.Item(41, 5).Value = "  -3.06%  "

$ws.Cells.Item(42, 4).Value = "3.36"
$ws.Cells.Item(42, 5).Value = "  -6.05%  "

$ws.Cells.Item(43, 4).Value = "266.74"
$ws.Cells.Item(43, 5).Value = "  -5.45%  "

$ws.Cells.Item(44, 4).Value = "4.92"
$ws.Cells.Item(44, 5).Value = "  -1.87%  "

$ws.Cells.Item(45, 4).Value = "0.582"
$ws.Cells.Item(45, 5).Value = "  -3.83%  "

$ws.Cells.Item(46, 4).Value = "122.09"
$ws.Cells.Item(46, 5).Value = "  -5.76%  "

$ws.Cells.Item(47, 4).Value = "0.0900"
$ws.Cells.Item(47, 5).Value = "  -2.71%  "

$ws.Cells.Item(48, 4).Value = "0.0481"
$ws.Cells.Item(48, 5).Value = "  -4.15%  "

$ws.Cells.Item(49, 4).Value = "0.0210"
$ws.Cells.Item(49, 5).Value = "  -3.53%  "

$ws.Cells.Item(50, 4).Value = "16.63"
$ws.Cells.Item(50, 5).Value = "  -3.90%  "

$ws.Cells.Item(51, 4).Value = "1.704.84"
$ws.Cells.Item(51, 5).Value = "  -3.19%  "

# Restore the original (default/"Normal") cell style now that the text
# is committed, so no stray number-format style lingers on the cells.
$dataRange.Style = "Normal"